$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.846.18"
$ws.Range("E2").Value = "  -0.47%  "
$ws.Range("D3").Value = "2.304.78"
$ws.Range("E3").Value = "  +0.21%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'307.01"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.29%  "
$ws.Range("D6").Value = "'96.52"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.44%  "
$ws.Range("E7").Value = "  -1.86%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("E9").Value = "  -2.32%  "
$ws.Range("D10").Value = "'35.45"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.84%  "
$ws.Range("D11").Value = "'0.0792"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.12%  "
$ws.Range("D12").Value = "'18.51"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +4.28%  "
$ws.Range("E13").Value = "  +1.37%  "
$ws.Range("E14").Value = "  -1.38%  "
$ws.Range("D15").Value = "2.663.48"
$ws.Range("E15").Value = "  +0.09%  "
$ws.Range("D16").Value = "2.307.11"
$ws.Range("E16").Value = "  +0.42%  "
$ws.Range("E17").Value = "  -0.41%  "
$ws.Range("D18").Value = "42.783.13"
$ws.Range("E18").Value = "  -0.41%  "
$ws.Range("D19").Value = "'13.16"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.83%  "
$ws.Range("D20").Value = "0.0₃0899"
$ws.Range("E20").Value = "  -1.17%  "
$ws.Range("E21").Value = "  -1.54%  "
$ws.Range("D22").Value = "'67.38"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.97%  "
$ws.Range("D23").Value = "'236.29"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.65%  "
$ws.Range("D24").Value = "'2.14"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.00%  "
$ws.Range("D25").Value = "'2.46"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.84%  "
$ws.Range("E26").Value = "  +0.07%  "
$ws.Range("D27").Value = "'4.00"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.40%  "
$ws.Range("D28").Value = "'25.27"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.21%  "
$ws.Range("D29").Value = "'2.38"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +16.49%  "
$ws.Range("D30").Value = "'166.21"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.97%  "
$ws.Range("E31").Value = "  -0.55%  "
$ws.Range("D32").Value = "'33.16"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.44%  "
$ws.Range("E33").Value = "  +0.03%  "
$ws.Range("E34").Value = "  -0.70%  "
$ws.Range("E35").Value = "  -2.25%  "
$ws.Range("D36").Value = "'17.79"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.15%  "
$ws.Range("E37").Value = "  -0.79%  "
$ws.Range("D38").Value = "'0.0694"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.55%  "
$ws.Range("E39").Value = "  -1.30%  "
$ws.Range("E40").Value = "  -0.88%  "
$ws.Range("E41").Value = "  -0.83%  "
$ws.Range("E42").Value = "  -2.65%  "
$ws.Range("D43").Value = "2.011.08"
$ws.Range("E43").Value = "  -0.30%  "
$ws.Range("E44").Value = "  -2.31%  "
$ws.Range("D45").Value = "'18.30"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +5.02%  "
$ws.Range("D46").Value = "'10.05"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.46%  "
$ws.Range("D47").Value = "'2.05"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -7.10%  "
$ws.Range("E48").Value = "  -0.77%  "
$ws.Range("E49").Value = "  +11.10%  "
$ws.Range("D50").Value = "'53.89"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.39%  "
$ws.Range("D51").Value = "2.528.12"
$ws.Range("E51").Value = "  -0.06%  "
